$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 70198.2
$ws.Range("J108").Value = 70198.2
$ws.Range("L108").Value = 70198.2
$ws.Range("N108").Value = -77878.2

$ws.Range("H112").Value = 2219.4
$ws.Range("J112").Value = 2219.4
$ws.Range("L112").Value = 6658.200000000001
$ws.Range("N112").Value = -8874.200000000001

$ws.Range("H117").Value = 78999
$ws.Range("J117").Value = 78999
$ws.Range("L117").Value = 78999
$ws.Range("N117").Value = -88177

$ws.Range("H133").Value = 105995.5
$ws.Range("J133").Value = 105995.5
$ws.Range("L133").Value = 105995.5
$ws.Range("N133").Value = -116115.5

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws.Range("H138").Value = 287794.47
$ws.Range("I138").Value = 3530.4688
$ws.Range("K138").Value = 10591.4064
$ws.Range("M138").Value = -5451.4064

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 684.6
$ws.Range("J2").Value = 797
$ws.Range("L2").Value = 797
$ws.Range("N2").Value = -1023

$ws.Range("H45").Value = 64048.57
$ws.Range("I45").Value = 85670
$ws.Range("J45").Value = 9995
$ws.Range("K45").Value = 85670
$ws.Range("L45").Value = 9995
$ws.Range("M45").Value = -85293
$ws.Range("N45").Value = -10749

$ws.Range("H63").Value = 3332.7778
$ws.Range("I63").Value = 2999.375
$ws.Range("K63").Value = 2999.375
$ws.Range("M63").Value = -2313.375

$ws.Range("H66").Value = 3332.7778
$ws.Range("I66").Value = 2999.375
$ws.Range("K66").Value = 14996.875
$ws.Range("M66").Value = -11564.875

$ws.Range("H97").Value = 1019.1739
$ws.Range("I97").Value = 1129.2778
$ws.Range("J97").Value = 622.8
$ws.Range("K97").Value = 1129.2778
$ws.Range("L97").Value = 622.8
$ws.Range("M97").Value = -633.2778000000001
$ws.Range("N97").Value = -1614.8

$ws.Range("H116").Value = 684.6
$ws.Range("J116").Value = 797
$ws.Range("L116").Value = 797
$ws.Range("N116").Value = -5385

$ws.Range("H132").Value = 2624.9092
$ws.Range("I132").Value = 1848.2858
$ws.Range("J132").Value = 3984
$ws.Range("K132").Value = 5544.857400000001
$ws.Range("L132").Value = 11952
$ws.Range("M132").Value = -3014.857400000001
$ws.Range("N132").Value = -17012

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 684.6
$ws.Range("J3").Value = 797
$ws.Range("L3").Value = 797
$ws.Range("N3").Value = -1025

$ws.Range("H53").Value = 80808
$ws.Range("J53").Value = 80808
$ws.Range("L53").Value = 80808
$ws.Range("N53").Value = -81956

$ws.Range("H86").Value = 3897.2593
$ws.Range("I86").Value = 3675.7
$ws.Range("J86").Value = 4530.2856
$ws.Range("K86").Value = 3675.7
$ws.Range("L86").Value = 4530.2856
$ws.Range("M86").Value = -2552.7
$ws.Range("N86").Value = -6776.2856

$ws.Range("H89").Value = 3897.2593
$ws.Range("I89").Value = 3675.7
$ws.Range("J89").Value = 4530.2856
$ws.Range("K89").Value = 18378.5
$ws.Range("L89").Value = 22651.428
$ws.Range("M89").Value = -12762.5
$ws.Range("N89").Value = -33883.428

$ws.Range("H94").Value = 71429290
$ws.Range("I94").Value = 71429290
$ws.Range("K94").Value = 71429290
$ws.Range("M94").Value = -71428839

$ws.Range("H105").Value = 8127352.5
$ws.Range("I105").Value = 436862.2
$ws.Range("K105").Value = 436862.2
$ws.Range("M105").Value = -435115.2

$ws.Range("H107").Value = 3664616.5
$ws.Range("I107").Value = 4809207
$ws.Range("J107").Value = 1926.6
$ws.Range("K107").Value = 4809207
$ws.Range("L107").Value = 1926.6
$ws.Range("M107").Value = -4807287
$ws.Range("N107").Value = -5766.6

$ws.Range("H134").Value = 3174.7856
$ws.Range("I134").Value = 2446.6191
$ws.Range("K134").Value = 7339.8573
$ws.Range("M134").Value = -4804.8573

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1807.2
$ws.Range("I16").Value = 1807.2
$ws.Range("K16").Value = 1807.2
$ws.Range("M16").Value = -1520.2

$ws.Range("H31").Value = 3739.1406
$ws.Range("I31").Value = 3218.1853
$ws.Range("J31").Value = 4119.2974
$ws.Range("K31").Value = 3218.1853
$ws.Range("L31").Value = 4119.2974
$ws.Range("M31").Value = -2923.1853
$ws.Range("N31").Value = -4709.2974

$ws.Range("H34").Value = 3739.1406
$ws.Range("I34").Value = 3218.1853
$ws.Range("J34").Value = 4119.2974
$ws.Range("K34").Value = 3218.1853
$ws.Range("L34").Value = 4119.2974
$ws.Range("M34").Value = -3016.1853
$ws.Range("N34").Value = -4523.2974

$ws.Range("H58").Value = 2981.111
$ws.Range("J58").Value = 3146.7646
$ws.Range("L58").Value = 3146.7646
$ws.Range("N58").Value = -3552.7646

$ws.Range("H62").Value = 14306893
$ws.Range("J62").Value = 30083.334
$ws.Range("L62").Value = 30083.334
$ws.Range("N62").Value = -31331.334

$ws.Range("H65").Value = 14306893
$ws.Range("J65").Value = 30083.334
$ws.Range("L65").Value = 150416.67
$ws.Range("N65").Value = -156656.67

$ws.Range("H99").Value = 6281.143
$ws.Range("I99").Value = 6328
$ws.Range("J99").Value = 6000
$ws.Range("K99").Value = 6328
$ws.Range("L99").Value = 6000
$ws.Range("M99").Value = -4830
$ws.Range("N99").Value = -8996

$ws.Range("H105").Value = 2545
$ws.Range("I105").Value = 1100
$ws.Range("K105").Value = 1100
$ws.Range("M105").Value = 647

$ws.Range("H113").Value = 1807.2
$ws.Range("I113").Value = 1807.2
$ws.Range("K113").Value = 1807.2
$ws.Range("M113").Value = 362.8

$ws.Range("H126").Value = 6281.143
$ws.Range("I126").Value = 6328
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 18984
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -16514
$ws.Range("N126").Value = -22940

$ws.Range("H132").Value = 11114284
$ws.Range("I132").Value = 13516346
$ws.Range("K132").Value = 40549038
$ws.Range("M132").Value = -40546508

$ws.Range("H136").Value = 2981.111
$ws.Range("J136").Value = 3146.7646
$ws.Range("L136").Value = 9440.293799999999
$ws.Range("N136").Value = -14540.2938

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 11804.99
$ws.Range("I4").Value = 12154.761
$ws.Range("J4").Value = 7208
$ws.Range("K4").Value = 36464.283
$ws.Range("L4").Value = 21624
$ws.Range("M4").Value = -36352.283
$ws.Range("N4").Value = -21848

$ws.Range("H52").Value = 1068.25
$ws.Range("J52").Value = 1068.25
$ws.Range("L52").Value = 3204.75
$ws.Range("N52").Value = -3736.75

$ws.Range("H132").Value = 2478.875
$ws.Range("J132").Value = 3187.4443
$ws.Range("L132").Value = 28686.9987
$ws.Range("N132").Value = -33746.9987

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 19998
$ws.Range("J5").Value = 19998
$ws.Range("L5").Value = 19998
$ws.Range("N5").Value = -20222

$ws.Range("H122").Value = 4531812.5
$ws.Range("I122").Value = 10992723
$ws.Range("J122").Value = 9174.799999999999
$ws.Range("K122").Value = 32978169
$ws.Range("L122").Value = 27524.4
$ws.Range("M122").Value = -32975719
$ws.Range("N122").Value = -32424.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 14344444
$ws.Range("J2").Value = 15762500
$ws.Range("L2").Value = 15762500
$ws.Range("N2").Value = -15762724

$ws.Range("H46").Value = 1858
$ws.Range("I46").Value = 3000
$ws.Range("K46").Value = 3000
$ws.Range("M46").Value = -2812

$ws.Range("H61").Value = 1429.2646
$ws.Range("I61").Value = 1369.24
$ws.Range("J61").Value = 1596
$ws.Range("K61").Value = 1369.24
$ws.Range("L61").Value = 1596
$ws.Range("M61").Value = -1167.24
$ws.Range("N61").Value = -2000

$ws.Range("H93").Value = 1963.5555
$ws.Range("I93").Value = 2253.1428
$ws.Range("J93").Value = 950
$ws.Range("K93").Value = 2253.1428
$ws.Range("L93").Value = 950
$ws.Range("M93").Value = -1005.1428
$ws.Range("N93").Value = -3446

$ws.Range("H113").Value = 1429.2646
$ws.Range("I113").Value = 1369.24
$ws.Range("J113").Value = 1596
$ws.Range("K113").Value = 1369.24
$ws.Range("L113").Value = 1596
$ws.Range("M113").Value = 800.76
$ws.Range("N113").Value = -5936

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 667866.7
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws.Range("H107").Value = 869.9091
$ws.Range("J107").Value = 529.5
$ws.Range("L107").Value = 1588.5
$ws.Range("N107").Value = -5428.5

$ws.Range("H132").Value = 8774972
$ws.Range("I132").Value = 10419864
$ws.Range("K132").Value = 31259592
$ws.Range("M132").Value = -31257062

$ws.Range("H135").Value = 53695.445
$ws.Range("J135").Value = 53695.445
$ws.Range("L135").Value = 53695.445
$ws.Range("N135").Value = -63835.445
